$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values rotate: C1<-"prediction", D1<-"rejection-f", E1<-"max"
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2 values: C2 becomes the text, D2 stays text, E2 becomes numeric 1
$ws.Range("C2").Value = "s__Massilicoli timonensis"
$ws.Range("D2").Value = "s__Massilicoli timonensis"
$ws.Range("E2").Value = 1
